function RGB($r,$g,$b) { return $r + ($g * 256) + ($b * 65536) }
$p = $ppt.ActivePresentation
$s1 = $p.Slides.Item(1)
$tcs = $s1.ThemeColorScheme
$c = $tcs.Colors(5)
Write-Output ("Before RGB: " + $c.RGB)
$c.RGB = (RGB 0x99 0xCB 0x38)
Write-Output ("After RGB: " + $c.RGB)
